# Generate Report for Handback
#
# This marks the two handed-off files as handed back in sync with en-US:
#  - updates the status text used on the Overview sheet,
#  - records the "Latest Target File" (source file) and "Latest Handback
#    File" (generated xlf) links/names for each localized file row on the
#    zh-cn and de-de sheets,
#  - stamps the handback datetime,
#  - widens a few columns so the longer strings are readable.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e7a089741d7d24873ec36150890917ba91f0b10d/e2e/139be023-839b-4dc9-b15a-fb184823e458.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e7a089741d7d24873ec36150890917ba91f0b10d/e2e/b39be596-6d61-4722-8a47-77695c6a1944.md"
$mdName1 = "139be023-839b-4dc9-b15a-fb184823e458.md"
$mdName2 = "b39be596-6d61-4722-8a47-77695c6a1944.md"

# --- Overview sheet: refresh the handback status shown for both locales ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# Status / Latest Target / Latest Handback columns read better wider now
# that they hold the new, longer text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Columns.Item(3).ColumnWidth = 29.14
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl1, $null, $null, $mdName1)
$wsZh.Range("J2").Value = "139be023-839b-4dc9-b15a-fb184823e458.308cd4b5d00aea98b892a0399a5b0d07830e0f0c.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-01 22:52:35"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl2, $null, $null, $mdName2)
$wsZh.Range("J3").Value = "b39be596-6d61-4722-8a47-77695c6a1944.22ebc7128498d347e65ea93b6119ca37db601d92.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-01 22:52:35"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(3).ColumnWidth = 29.14
$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl1, $null, $null, $mdName1)
$wsDe.Range("J2").Value = "139be023-839b-4dc9-b15a-fb184823e458.308cd4b5d00aea98b892a0399a5b0d07830e0f0c.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-01 22:52:41"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl2, $null, $null, $mdName2)
$wsDe.Range("J3").Value = "b39be596-6d61-4722-8a47-77695c6a1944.22ebc7128498d347e65ea93b6119ca37db601d92.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-01 22:52:41"
